$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1722.4375
$ws.Range("I41").Value = 411.6
$ws.Range("J41").Value = 2318.2727
$ws.Range("K41").Value = 411.6
$ws.Range("L41").Value = 2318.2727
$ws.Range("M41").Value = 28.39999999999998
$ws.Range("N41").Value = -3198.2727
$ws.Range("H137").Value = 11841031
$ws.Range("I137").Value = 1001071.8
$ws.Range("K137").Value = 3003215.4
$ws.Range("M137").Value = -3000665.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2707.69
$ws.Range("I32").Value = 2500.7917
$ws.Range("K32").Value = 2500.7917
$ws.Range("M32").Value = -2213.7917
$ws.Range("H45").Value = 3159.4827
$ws.Range("I45").Value = 2985.8
$ws.Range("K45").Value = 2985.8
$ws.Range("M45").Value = -2608.8
$ws.Range("H63").Value = 400
$ws.Range("J63").Value = 200
$ws.Range("L63").Value = 200
$ws.Range("N63").Value = -1572
$ws.Range("H66").Value = 400
$ws.Range("J66").Value = 200
$ws.Range("L66").Value = 1000
$ws.Range("N66").Value = -7864
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 67987.5
$ws.Range("J80").Value = 67987.5
$ws.Range("L80").Value = 67987.5
$ws.Range("N80").Value = -69983.5
$ws.Range("H82").Value = 25000
$ws.Range("I82").Value = 25000
$ws.Range("K82").Value = 25000
$ws.Range("M82").Value = -24639
$ws.Range("H83").Value = 67987.5
$ws.Range("J83").Value = 67987.5
$ws.Range("L83").Value = 203962.5
$ws.Range("N83").Value = -213946.5
$ws.Range("H85").Value = 25000
$ws.Range("I85").Value = 25000
$ws.Range("K85").Value = 25000
$ws.Range("M85").Value = -23752
$ws.Range("H88").Value = 34697.625
$ws.Range("J88").Value = 36714.535
$ws.Range("L88").Value = 36714.535
$ws.Range("N88").Value = -37526.535
$ws.Range("H91").Value = 34697.625
$ws.Range("J91").Value = 36714.535
$ws.Range("L91").Value = 36714.535
$ws.Range("N91").Value = -39522.535
$ws.Range("H97").Value = 1065.5714
$ws.Range("I97").Value = 1093.1666
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 1093.1666
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -597.1666
$ws.Range("N97").Value = -1892
$ws.Range("H132").Value = 3421.1853
$ws.Range("I132").Value = 1772.5454
$ws.Range("K132").Value = 5317.6362
$ws.Range("M132").Value = -2787.6362
$ws.Range("H135").Value = 98999.336
$ws.Range("J135").Value = 98999.336
$ws.Range("L135").Value = 98999.336
$ws.Range("N135").Value = -109139.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3000.9773
$ws.Range("I20").Value = 3154.5
$ws.Range("J20").Value = 2894.6924
$ws.Range("K20").Value = 3154.5
$ws.Range("L20").Value = 2894.6924
$ws.Range("M20").Value = -2907.5
$ws.Range("N20").Value = -3388.6924
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H82").Value = 12279
$ws.Range("I82").Value = 12279
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 12279
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -11896
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 12279
$ws.Range("I85").Value = 12279
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 12279
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -10953
$ws.Range("N85").ClearContents()
$ws.Range("H134").Value = 2581.818
$ws.Range("I134").Value = 2538.1191
$ws.Range("K134").Value = 7614.3573
$ws.Range("M134").Value = -5079.3573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 315.92307
$ws.Range("I7").Value = 540.4
$ws.Range("J7").Value = 175.625
$ws.Range("K7").Value = 540.4
$ws.Range("L7").Value = 175.625
$ws.Range("M7").Value = -427.4
$ws.Range("N7").Value = -401.625
$ws.Range("H132").Value = 26145806
$ws.Range("I132").Value = 28987212
$ws.Range("J132").Value = 4879.4
$ws.Range("K132").Value = 86961636
$ws.Range("L132").Value = 14638.2
$ws.Range("M132").Value = -86959106
$ws.Range("N132").Value = -19698.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 36239
$ws.Range("J101").Value = 36239
$ws.Range("L101").Value = 108717
$ws.Range("N101").Value = -113585
$ws.Range("H131").Value = 13597220
$ws.Range("I131").Value = 6668246
$ws.Range("J131").Value = 17625694
$ws.Range("K131").Value = 20004738
$ws.Range("L131").Value = 52877082
$ws.Range("M131").Value = -19999698
$ws.Range("N131").Value = -52887162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 416.44446
$ws.Range("I2").Value = 220.4
$ws.Range("J2").Value = 661.5
$ws.Range("K2").Value = 220.4
$ws.Range("L2").Value = 661.5
$ws.Range("M2").Value = -107.4
$ws.Range("N2").Value = -887.5
$ws.Range("H57").Value = 20998.4
$ws.Range("J57").Value = 20998.4
$ws.Range("L57").Value = 20998.4
$ws.Range("N57").Value = -22638.4
$ws.Range("H80").Value = 76116.57000000001
$ws.Range("I80").Value = 104378.7
$ws.Range("K80").Value = 104378.7
$ws.Range("M80").Value = -103380.7
$ws.Range("H83").Value = 76116.57000000001
$ws.Range("I83").Value = 104378.7
$ws.Range("K83").Value = 521893.5
$ws.Range("M83").Value = -516901.5
$ws.Range("H96").Value = 21261
$ws.Range("J96").Value = 21261
$ws.Range("L96").Value = 21261
$ws.Range("N96").Value = -26753
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H132").Value = 113466.11
$ws.Range("I132").Value = 201097.1
$ws.Range("J132").Value = 3927.375
$ws.Range("K132").Value = 603291.3
$ws.Range("L132").Value = 11782.125
$ws.Range("M132").Value = -600761.3
$ws.Range("N132").Value = -16842.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 69696
$ws.Range("J130").Value = 69696
$ws.Range("L130").Value = 69696
$ws.Range("N130").Value = -79736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 39000
$ws.Range("J54").Value = 38000
$ws.Range("L54").Value = 38000
$ws.Range("N54").Value = -39040
$ws.Range("H122").Value = 3886.5151
$ws.Range("I122").Value = 3645.6316
$ws.Range("J122").Value = 4213.4287
$ws.Range("K122").Value = 10936.8948
$ws.Range("L122").Value = 12640.2861
$ws.Range("M122").Value = -8486.8948
$ws.Range("N122").Value = -17540.2861
$ws.Range("H128").Value = 119999.75
$ws.Range("J128").Value = 119999.75
$ws.Range("L128").Value = 119999.75
$ws.Range("N128").Value = -129959.75
$ws.Range("H136").Value = 3177.9187
$ws.Range("I136").Value = 2033.1167
$ws.Range("J136").Value = 5819.769
$ws.Range("K136").Value = 6099.3501
$ws.Range("L136").Value = 17459.307
$ws.Range("M136").Value = -3549.3501
$ws.Range("N136").Value = -22559.307
